$wb = $excel.ActiveWorkbook

# Edits extracted from the canonical OOXML diff: cell -> new value, per sheet.
# Format: @("<CellRef>", <NewValue>)
$edits = @{}

$edits["ALC"] = @(
    @("H15", 18519742),
    @("I15", 18519742),
    @("K15", 55559226),
    @("M15", -55559057),
    @("H33", 1630.8572),
    @("I33", 1713),
    @("J33", 1138),
    @("K33", 1713),
    @("L33", 1138),
    @("M33", -1484),
    @("N33", -1596),
    @("H40", 149057),
    @("I40", 5500),
    @("J40", 172983.17),
    @("K40", 5500),
    @("L40", 172983.17),
    @("M40", -5325),
    @("N40", -173333.17),
    @("H132", 1235.125),
    @("I132", 1042.6),
    @("J132", 2197.75),
    @("K132", 3127.8),
    @("L132", 6593.25),
    @("M132", -597.7999999999997),
    @("N132", -11653.25),
    @("H137", 3834.8823),
    @("J137", 6999.6665),
    @("L137", 20998.9995),
    @("N137", -26098.9995)
)

$edits["ARM"] = @(
    @("H61", 33336024),
    @("I61", 1264.091),
    @("J61", 125006616),
    @("K61", 1264.091),
    @("L61", 125006616),
    @("M61", -1052.091),
    @("N61", -125007040),
    @("H74", 22930.041),
    @("I74", 27026.025),
    @("J74", 5180.778),
    @("K74", 27026.025),
    @("L74", 5180.778),
    @("M74", -26152.025),
    @("N74", -6928.778),
    @("H77", 22930.041),
    @("I77", 27026.025),
    @("J77", 5180.778),
    @("K77", 135130.125),
    @("L77", 25903.89),
    @("M77", -130762.125),
    @("N77", -34639.89),
    @("H102", 3032.6956),
    @("I102", 2587.6667),
    @("K102", 2587.6667),
    @("M102", -965.6667000000002),
    @("H122", 3384.8809),
    @("I122", 2068.5),
    @("K122", 6205.5),
    @("M122", -3755.5),
    @("H132", 3723.7058),
    @("I132", 2690.9636),
    @("K132", 8072.8908),
    @("M132", -5542.8908),
    @("H136", 33336024),
    @("I136", 1264.091),
    @("J136", 125006616),
    @("K136", 3792.273),
    @("L136", 375019848),
    @("M136", -1242.273),
    @("N136", -375024948)
)

$edits["BSM"] = @(
    @("H94", 1472.4375),
    @("I94", 639.14703),
    @("K94", 639.14703),
    @("M94", -188.14703),
    @("H107", 35166890),
    @("I107", 51146564),
    @("J107", 11608),
    @("K107", 51146564),
    @("L107", 11608),
    @("M107", -51144644),
    @("N107", -15448),
    @("H134", 6414392),
    @("I134", 11365440),
    @("K134", 34096320),
    @("M134", -34093785)
)

$edits["CRP"] = @(
    @("H58", 4296.727),
    @("J58", 5769),
    @("L58", 5769),
    @("N58", -6175),
    @("H59", 55664),
    @("J59", 55664),
    @("L59", 55664),
    @("N59", -57954),
    @("H76", 5172),
    @("I76", 5172),
    @("K76", 5172),
    @("M76", -4857),
    @("H79", 5172),
    @("I79", 5172),
    @("K79", 5172),
    @("M79", -4080),
    @("H93", 7133.3335),
    @("I93", 3200),
    @("J93", 15000),
    @("K93", 3200),
    @("L93", 15000),
    @("M93", -1328),
    @("N93", -18744),
    @("H132", 2984.02),
    @("I132", 2155.3513),
    @("K132", 6466.053899999999),
    @("M132", -3936.053899999999),
    @("H134", 3494.2703),
    @("I134", 2010.4642),
    @("K134", 6031.392599999999),
    @("M134", -3496.392599999999),
    @("H136", 4296.727),
    @("J136", 5769),
    @("L136", 17307),
    @("N136", -22407)
)

$edits["CUL"] = @(
    @("H80", 21982608),
    @("I80", 347687.5),
    @("J80", 71433860),
    @("K80", 1043062.5),
    @("L80", 214301580),
    @("M80", -1042126.5),
    @("N80", -214303452),
    @("H83", 21982608),
    @("I83", 347687.5),
    @("J83", 71433860),
    @("K83", 3129187.5),
    @("L83", 642904740),
    @("M83", -3124507.5),
    @("N83", -642914100),
    @("H86", 828.5),
    @("I86", 361),
    @("J86", 984.3333),
    @("K86", 1083),
    @("L86", 2952.9999),
    @("M86", 103),
    @("N86", -5324.9999),
    @("H89", 828.5),
    @("I89", 361),
    @("J89", 984.3333),
    @("K89", 3249),
    @("L89", 8858.9997),
    @("M89", 2679),
    @("N89", -20714.9997),
    @("H93", 6981),
    @("H124", 3465.2856),
    @("I124", 3465.2856),
    @("K124", 10395.8568),
    @("M124", -5485.856800000001),
    @("H125", 4400),
    @("J125", 5800),
    @("L125", 17400),
    @("N125", -27240),
    @("H132", 5193.391),
    @("I132", 2519.077),
    @("J132", 8670),
    @("K132", 22671.693),
    @("L132", 78030),
    @("M132", -20141.693),
    @("N132", -83090),
    @("H138", 71151.664),
    @("I138", 71151.664),
    @("K138", 213454.992),
    @("M138", -208314.992),
    @("H139", 51471.094),
    @("J139", 9999.4),
    @("L139", 29998.2),
    @("N139", -40278.2),
    @("H140", 113195.445),
    @("J140", 6332),
    @("L140", 18996),
    @("N140", -29356),
    @("H141", 2896.5),
    @("I141", 2896.5),
    @("K141", 8689.5),
    @("M141", -3509.5)
)

$edits["GSM"] = @(
    @("H75", 30929.857),
    @("J75", 30929.857),
    @("L75", 30929.857),
    @("N75", -32677.857),
    @("H78", 30929.857),
    @("J78", 30929.857),
    @("L78", 92789.571),
    @("N78", -101525.571),
    @("H122", 1611989.5),
    @("I122", 1959829.5),
    @("K122", 5879488.5),
    @("M122", -5877038.5),
    @("H132", 2482.15),
    @("I132", 1488.0769),
    @("J132", 4328.2856),
    @("K132", 4464.2307),
    @("L132", 12984.8568),
    @("M132", -1934.2307),
    @("N132", -18044.8568)
)

$edits["LTW"] = @(
    @("H46", 2641.3333),
    @("I46", 2099.889),
    @("K46", 2099.889),
    @("M46", -1911.889),
    @("H69", 53961),
    @("J69", 53961),
    @("L69", 53961),
    @("N69", -55583),
    @("H72", 53961),
    @("J72", 53961),
    @("L72", 161883),
    @("N72", -169995),
    @("H93", 1188.2424),
    @("I93", 1193.1364),
    @("K93", 1193.1364),
    @("M93", 54.86359999999991),
    @("H100", 3581.8333),
    @("I100", 2246.9412),
    @("J100", 6823.7144),
    @("K100", 2246.9412),
    @("L100", 6823.7144),
    @("M100", -1705.9412),
    @("N100", -7905.7144),
    @("H107", 3247.7778),
    @("I107", 3247.7778),
    @("K107", 3247.7778),
    @("M107", -1327.7778),
    @("H132", 10425137),
    @("J132", 12597.407),
    @("L132", 37792.221),
    @("N132", -42852.221)
)

$edits["WVR"] = @(
    @("H126", 200009630),
    @("I126", 333340800),
    @("K126", 1000022400),
    @("M126", -1000019930)
)

foreach ($sheetName in $edits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $edits[$sheetName]) {
        $cellRef = $pair[0]
        $newValue = $pair[1]
        $ws.Range($cellRef).Value = $newValue
    }
}

Write-Output "Applied $($edits.Values.Count) sheet groups of edits."